# Apply weekly fruit/vegetable price updates (Hortaliza, Femacal de La Calera - Zanahoria)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 410
$ws.Range("D410").Value = 44736
$ws.Range("J410").Value = 480
$ws.Range("K410").Value = 7000
$ws.Range("L410").Value = 7300
$ws.Range("M410").Value = 7156
$ws.Range("P410").Value = 358

# Row 411
$ws.Range("D411").Value = 44299
$ws.Range("J411").Value = 300
$ws.Range("K411").Value = 6000
$ws.Range("L411").Value = 6500
$ws.Range("M411").Value = 6417
$ws.Range("P411").Value = 321

# Row 412
$ws.Range("D412").Value = 44334
$ws.Range("J412").Value = 310
$ws.Range("K412").Value = 5000
$ws.Range("L412").Value = 5500
$ws.Range("M412").Value = 5242
$ws.Range("O412").Value = "Chillán"
$ws.Range("P412").Value = 262

# Row 413
$ws.Range("D413").Value = 45119
$ws.Range("J413").Value = 280
$ws.Range("K413").Value = 7000
$ws.Range("L413").Value = 7500
$ws.Range("M413").Value = 7214
$ws.Range("O413").Value = "Calera"
$ws.Range("P413").Value = 361

# Row 414
$ws.Range("D414").Value = 44582
$ws.Range("J414").Value = 760
$ws.Range("M414").Value = 7250
$ws.Range("O414").Value = "Provincia de Quillota"
$ws.Range("P414").Value = 362

# Row 415
$ws.Range("D415").Value = 44263
$ws.Range("J415").Value = 530
$ws.Range("K415").Value = 6000
$ws.Range("L415").Value = 6500
$ws.Range("M415").Value = 6264
$ws.Range("P415").Value = 313

# Row 416
$ws.Range("D416").Value = 44210
$ws.Range("J416").Value = 250
$ws.Range("K416").Value = 8000
$ws.Range("L416").Value = 8500
$ws.Range("M416").Value = 8260
$ws.Range("P416").Value = 413

# Row 417
$ws.Range("D417").Value = 44921
$ws.Range("J417").Value = 310
$ws.Range("K417").Value = 10500
$ws.Range("L417").Value = 11000
$ws.Range("M417").Value = 10758
$ws.Range("P417").Value = 538

# Row 418
$ws.Range("D418").Value = 44907
$ws.Range("J418").Value = 510
$ws.Range("K418").Value = 11000
$ws.Range("L418").Value = 12000
$ws.Range("M418").Value = 11510
$ws.Range("P418").Value = 576

# Row 419
$ws.Range("D419").Value = 44550
$ws.Range("J419").Value = 310
$ws.Range("K419").Value = 6500
$ws.Range("L419").Value = 7000
$ws.Range("M419").Value = 6742
$ws.Range("P419").Value = 337

# Row 420
$ws.Range("D420").Value = 44813
$ws.Range("J420").Value = 240
$ws.Range("K420").Value = 11500
$ws.Range("L420").Value = 12000
$ws.Range("M420").Value = 11750
$ws.Range("P420").Value = 588

# Row 421
$ws.Range("D421").Value = 44483
$ws.Range("J421").Value = 510
$ws.Range("K421").Value = 8000
$ws.Range("L421").Value = 8500
$ws.Range("M421").Value = 8255
$ws.Range("O421").Value = "Chillán"
$ws.Range("P421").Value = 413

# Row 422
$ws.Range("D422").Value = 44875
$ws.Range("J422").Value = 180
$ws.Range("K422").Value = 13000
$ws.Range("L422").Value = 13000
$ws.Range("M422").Value = 13000
$ws.Range("O422").Value = "Provincia de Quillota"
$ws.Range("P422").Value = 650

# Row 423
$ws.Range("D423").Value = 44270
$ws.Range("J423").Value = 540
$ws.Range("K423").Value = 6000
$ws.Range("L423").Value = 6500
$ws.Range("M423").Value = 6269
$ws.Range("P423").Value = 313

# Row 424
$ws.Range("D424").Value = 44431
$ws.Range("J424").Value = 2083
$ws.Range("K424").Value = 4500
$ws.Range("L424").Value = 5000
$ws.Range("M424").Value = 4957
$ws.Range("P424").Value = 248

# Row 425
$ws.Range("D425").Value = 45049
$ws.Range("J425").Value = 310
$ws.Range("K425").Value = 8000
$ws.Range("L425").Value = 8500
$ws.Range("M425").Value = 8258
$ws.Range("P425").Value = 413

# Row 426
$ws.Range("D426").Value = 44264
$ws.Range("J426").Value = 250
$ws.Range("K426").Value = 6000
$ws.Range("L426").Value = 6000
$ws.Range("M426").Value = 6000
$ws.Range("P426").Value = 300

# Row 427
$ws.Range("D427").Value = 44967
$ws.Range("J427").Value = 240
$ws.Range("K427").Value = 8000
$ws.Range("L427").Value = 8500
$ws.Range("M427").Value = 8250
$ws.Range("P427").Value = 412

# Row 428
$ws.Range("D428").Value = 44306
$ws.Range("J428").Value = 160
$ws.Range("K428").Value = 5500
$ws.Range("L428").Value = 5500
$ws.Range("M428").Value = 5500
$ws.Range("P428").Value = 275

# Row 429
$ws.Range("D429").Value = 44516
$ws.Range("J429").Value = 540
$ws.Range("K429").Value = 6000
$ws.Range("L429").Value = 6500
$ws.Range("M429").Value = 6259
$ws.Range("P429").Value = 313

# Row 430
$ws.Range("D430").Value = 44991
$ws.Range("J430").Value = 290
$ws.Range("K430").Value = 7500
$ws.Range("L430").Value = 8000
$ws.Range("M430").Value = 7759
$ws.Range("P430").Value = 388

# Row 431
$ws.Range("D431").Value = 44435
$ws.Range("J431").Value = 3393
$ws.Range("K431").Value = 4500
$ws.Range("L431").Value = 5500
$ws.Range("M431").Value = 4979
$ws.Range("P431").Value = 249

# Row 432
$ws.Range("D432").Value = 44588
$ws.Range("J432").Value = 430
$ws.Range("K432").Value = 6500
$ws.Range("L432").Value = 7000
$ws.Range("M432").Value = 6733
$ws.Range("P432").Value = 337

# Row 433
$ws.Range("D433").Value = 44225
$ws.Range("J433").Value = 65
$ws.Range("L433").Value = 6500
$ws.Range("M433").Value = 6500
$ws.Range("P433").Value = 325

# Row 434
$ws.Range("D434").Value = 45120
$ws.Range("J434").Value = 230
$ws.Range("K434").Value = 7000
$ws.Range("L434").Value = 7500
$ws.Range("M434").Value = 7261
$ws.Range("O434").Value = "Calera"
$ws.Range("P434").Value = 363

# Row 582
$ws.Range("D582").Value = 45121
$ws.Range("J582").Value = 120
$ws.Range("K582").Value = 7000
$ws.Range("L582").Value = 7000
$ws.Range("M582").Value = 7000
$ws.Range("O582").Value = "Calera"
$ws.Range("P582").Value = 350

# Row 583
$ws.Range("D583").Value = 44341
$ws.Range("J583").Value = 180
$ws.Range("K583").Value = 5500
$ws.Range("L583").Value = 5500
$ws.Range("M583").Value = 5500
$ws.Range("P583").Value = 275

# Row 584
$ws.Range("D584").Value = 44777
$ws.Range("J584").Value = 340
$ws.Range("K584").Value = 11500
$ws.Range("L584").Value = 12000
$ws.Range("M584").Value = 11779
$ws.Range("P584").Value = 589

# Row 585
$ws.Range("D585").Value = 44662
$ws.Range("J585").Value = 310
$ws.Range("K585").Value = 6500
$ws.Range("L585").Value = 7000
$ws.Range("M585").Value = 6758
$ws.Range("P585").Value = 338

# Row 586
$ws.Range("D586").Value = 44607
$ws.Range("E586").Value = 5
$ws.Range("F586").Value = 100114013
$ws.Range("G586").Value = "Zanahoria"
$ws.Range("H586").Value = "Sin especificar"
$ws.Range("I586").Value = "Primera"
$ws.Range("J586").Value = 130
$ws.Range("K586").Value = 9000
$ws.Range("L586").Value = 10000
$ws.Range("M586").Value = 9538
$ws.Range("N586").Value = "`$/saco 20 kilos"
$ws.Range("O586").Value = "Provincia de Quillota"
$ws.Range("P586").Value = 477
$ws.Range("Q586").Value = 20
$ws.Range("R586").Value = "Hortaliza"

# Row 587
$ws.Range("A587").Value = 3
$ws.Range("B587").Value = "Femacal de La Calera"
$ws.Range("C587").Value = "Coquimbo"
$ws.Range("D587").Value = 45072
$ws.Range("E587").Value = 5
$ws.Range("F587").Value = 100114013
$ws.Range("G587").Value = "Zanahoria"
$ws.Range("H587").Value = "Sin especificar"
$ws.Range("I587").Value = "Primera"
$ws.Range("J587").Value = 230
$ws.Range("K587").Value = 7000
$ws.Range("L587").Value = 7500
$ws.Range("M587").Value = 7261
$ws.Range("N587").Value = "`$/saco 20 kilos"
$ws.Range("O587").Value = "Provincia de Quillota"
$ws.Range("P587").Value = 363
$ws.Range("Q587").Value = 20
$ws.Range("R587").Value = "Hortaliza"

# New row 587 needs the same date number format as the rest of column D
$ws.Range("D587").NumberFormat = $ws.Range("D586").NumberFormat
